$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Audio settings block (rows 8-15) -------------------------------------
# Row 8: header changed from "SETTINGS_audio"/"Audio" -> "SETTINGS_audio"/"Audio Settings"
$ws.Range("A8").Value = "SETTINGS_audio"
$ws.Range("B8").Value = "Audio Settings"

# Row 9: master volume
$ws.Range("A9").Value = "SETTINGS_audio.master"
$ws.Range("B9").Value = "Master Volume"

# Row 10: music volume
$ws.Range("A10").Value = "SETTINGS_audio.music"
$ws.Range("B10").Value = "Music Volume"

# Row 11: new menu sfx volume entry (replaces old sfx/players rows)
$ws.Range("A11").Value = "SETTINGS_audio.menu_sfx"
$ws.Range("B11").Value = "Menu SFX Volume"

# Row 12: new game sfx volume entry
$ws.Range("A12").Value = "SETTINGS_audio.game_sfx"
$ws.Range("B12").Value = "Game SFX Volume"

# Row 13: output device (moved up from row 17)
$ws.Range("A13").Value = "SETTINGS_audio.output"
$ws.Range("B13").Value = "Output Device"

# Row 14: reduce high frequencies (moved up from row 18)
$ws.Range("A14").Value = "SETTINGS_audio.reduceHigh"
$ws.Range("B14").Value = "Reduce high frequencies"

# Row 15: reduce high frequencies description (moved up from row 19)
$ws.Range("A15").Value = "SETTINGS_audio.reduceHigh.desc"
$ws.Range("B15").Value = "Applies a low-pass filter to all audio to reduce high pitched noises"

# Rows 16-19 (old input-device / ptt / open-mic rows) are fully removed.
$ws.Range("A16:B19").Clear()

# --- Settings category labels now get a value column too -----------------
$ws.Range("B21").Value = "Video"
$ws.Range("B23").Value = "Controls"
$ws.Range("B25").Value = "Accessibility"

# --- Dialogue test resources ------------------------------------------------
# Row 29 key renamed from DIALOGUE_test/conv1/1 -> DIALOGUE_test/conv1/d1
$ws.Range("A29").Value = "DIALOGUE_test/conv1/d1"
$ws.Range("B29").Value = "This is a test dialogue string"

# New row 30: second line of the test dialogue
$ws.Range("A30").Value = "DIALOGUE_test/conv1/d2"
$ws.Range("B30").Value = "What the freak"

# --- Update view / selection to match the edited file ----------------------
# Scroll so row 4 is at the top of the visible window, then select A15
# (best-effort; some hosts don't persist window scroll position to the XML).
$win = $excel.ActiveWindow
try {
    $win.ScrollRow = 4
    $win.ScrollColumn = 1
} catch {
}
$null = $ws.Range("A15").Select()
